# Actualizacion Datos Personales 4 nov
# Update row 2 (Camacho Juárez Sergio Eduardo / 5ARHV) statistics
# on both the "1er Parcial" and "3er Parcial" sheets.

$wb = $excel.ActiveWorkbook

$sheetNames = @("1er Parcial", "3er Parcial")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("E2").Value = 25
    $ws.Range("F2").Value = 12
    $ws.Range("G2").Value = 67.56999999999999
    $ws.Range("H2").Value = 32.43
    $ws.Range("I2").Value = 7.7
    $ws.Range("J2").Value = 0
    $ws.Range("K2").Value = 0
}
